# Insert a new row at position 68, shifting existing rows 68:183 down to 69:184.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("68").Insert()

# Populate the newly inserted row 68 with the new daily record.
$ws.Range("A68").Value = 5
$ws.Range("B68").Value = "Macroferia Regional de Talca"
$ws.Range("C68").Value = "Maule"
$ws.Range("D68").Value = 44477
$ws.Range("D68").NumberFormat = $ws.Range("D69").NumberFormat
$ws.Range("E68").Value = 7
$ws.Range("F68").Value = 100114014
$ws.Range("G68").Value = "Betarraga"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 650
$ws.Range("L68").Value = 650
$ws.Range("M68").Value = 650
$ws.Range("N68").Value = "$/paquete 5 unidades"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 130
$ws.Range("Q68").Value = 5
$ws.Range("R68").Value = "Hortaliza"
